$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 85: copy formatting (incl. column A) from row 83 pattern (s=4/5, no border) ---
$ws.Range("A83:E83").Copy()
$ws.Range("A85:E85").PasteSpecial(-4122)

# --- Row 86: copy formatting (B:E only, no column A cell) from row 80 pattern ---
$ws.Range("B80:E80").Copy()
$ws.Range("B86:E86").PasteSpecial(-4122)

# --- Rows 87-90: copy formatting (B:E only) from rows 81-84 pattern (reuse Yes/No/etc look) ---
$ws.Range("B81:E81").Copy()
$ws.Range("B87:E87").PasteSpecial(-4122)
$ws.Range("B82:E82").Copy()
$ws.Range("B88:E88").PasteSpecial(-4122)
$ws.Range("B83:E83").Copy()
$ws.Range("B89:E89").PasteSpecial(-4122)
$ws.Range("B80:E80").Copy()
$ws.Range("B90:E90").PasteSpecial(-4122)

# --- New shared strings are introduced in this precise order to reproduce the
#     original sharedStrings.xml ordering (index 235..241): C85,C86,D85,D86,E86,E85,A85 ---
$ws.Cells.Item(85, 3).Value = " HEY! It\'s you two![K] Came to visit\nyour old friends at the guild?"
$ws.Cells.Item(86, 3).Value = " WHAT?[K] You want to pull sentry\nduty and earn some money?"
$ws.Cells.Item(85, 4).Value = " ЭЙ! Это же вы, двое![K] Решили\nнавестить своих старых гильдейских\nтоварищей?"
$ws.Cells.Item(86, 4).Value = " ЧТО?[K] Хотите подзаработать,\nработая стражами?"
$ws.Cells.Item(86, 5).Value = " ŒÓÏ?[K] Öïóéóå ðïäèàñàáïóàóû,\nñàáïóàÿ òóñàçàíé?"
$ws.Cells.Item(85, 5).Value = " ÜÊ! Üóï çå âú, äâïå![K] Ñåšéìé\nîàâåòóéóû òâïéö òóàñúö ãéìûäåêòëéö\nóïâàñéþåê?"
$ws.Cells.Item(85, 1).Value = "SCRIPT/G01P04A/us2302.ssb"

# --- Remaining numeric cells (do not create new shared strings) ---
$ws.Cells.Item(85, 2).Value = 18
$ws.Cells.Item(86, 2).Value = 21

# --- Row 87 values (reuses existing shared strings Yes / Да / Äà) ---
$ws.Cells.Item(87, 2).Value = 25
$ws.Cells.Item(87, 3).Value = "Yes"
$ws.Cells.Item(87, 4).Value = "Да"
$ws.Cells.Item(87, 5).Value = "Äà"

# --- Row 88 values (reuses existing shared strings No / Нет / Îåó) ---
$ws.Cells.Item(88, 2).Value = 41
$ws.Cells.Item(88, 3).Value = "No"
$ws.Cells.Item(88, 4).Value = "Нет"
$ws.Cells.Item(88, 5).Value = "Îåó"

# --- Row 89 values (reuses existing shared strings "Oh, you won't, HUH?" set) ---
$ws.Cells.Item(89, 2).Value = 46
$ws.Cells.Item(89, 3).Value = " Oh, you won\'t, HUH?"
$ws.Cells.Item(89, 4).Value = " О, значит нет, А?"
$ws.Cells.Item(89, 5).Value = " Ï, èîàœéó îåó, À?"

# --- Row 90 values (reuses existing shared strings Check High Score / etc.) ---
$ws.Cells.Item(90, 2).Value = 51
$ws.Cells.Item(90, 3).Value = "Check High Score"
$ws.Cells.Item(90, 4).Value = "Таблица Рекордов"
$ws.Cells.Item(90, 5).Value = "Óàáìéøà Ñåëïñäïâ"

# --- Row heights to match autosized wrap-text content (engine does not auto-recompute for new rows) ---
$ws.Rows.Item(85).RowHeight = 43.2
$ws.Rows.Item(86).RowHeight = 21.6

# --- Update selection to match new active cell ---
$ws.Range("D85").Select()

Write-Host "Edit complete"
